# ---------------------------------------------------------------------------
# Daily Details - shift the schedule window from Aug 2025 to Oct 2025 and
# restyle the "continuation" date cells in column A (fill + left/center
# alignment, regular - not bold - weight) to match the first row of each
# day-group.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Formatting: the "continuation" rows for each day (everything except the
#    first row of a day-group) currently render bold-free already, but they
#    are missing the light-grey fill + left/center alignment that the first
#    row of each group uses. Bring them into line, and drop the leftover
#    bold weight from the day-group header font (A2/A7/A10/A14) so every date
#    cell in the column now shares the same regular weight.
# ---------------------------------------------------------------------------

$fillColor = 15132390   # 0x00E6E6E6 light grey, same as the header groups

# un-bold the first-of-day-group cells (font shared with the cells below)
$ws.Range("A2").Font.Bold = $false
$ws.Range("A7").Font.Bold = $false
$ws.Range("A10").Font.Bold = $false
$ws.Range("A14").Font.Bold = $false

# numFmtId 165 (mm/dd/yyyy) continuation cells, borderId 5
$r = $ws.Range("A3:A5")
$r.Interior.Color = $fillColor
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4108
$r.Font.Bold = $false

$r = $ws.Range("A8")
$r.Interior.Color = $fillColor
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4108
$r.Font.Bold = $false

$r = $ws.Range("A11:A13")
$r.Interior.Color = $fillColor
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4108
$r.Font.Bold = $false

# numFmtId 165 (mm/dd/yyyy) continuation cells, borderId 8
$r = $ws.Range("A6")
$r.Interior.Color = $fillColor
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4108
$r.Font.Bold = $false

$r = $ws.Range("A9")
$r.Interior.Color = $fillColor
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4108
$r.Font.Bold = $false

# numFmtId 164 (yyyy-mm-dd) continuation cell, borderId 5
$r = $ws.Range("A15")
$r.Interior.Color = $fillColor
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4108
$r.Font.Bold = $false

# numFmtId 164 (yyyy-mm-dd) continuation cell, borderId 8
$r = $ws.Range("A16")
$r.Interior.Color = $fillColor
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4108
$r.Font.Bold = $false

# ---------------------------------------------------------------------------
# 2) Data: push every date in the sheet forward by 60 days (2025-08-02 through
#    2025-08-04 become 2025-10-01 through 2025-10-03), updating both the
#    numeric Date column (A) and the "Unique Identifier" text column
#    (V for the first batch of rows, W for the rest) that embeds the same
#    date as its first pipe-delimited field.
# ---------------------------------------------------------------------------

function Shift-DailyDetailsRow {
    param(
        [int]$Row,
        [string]$IdCol
    )

    $dateCell = $ws.Cells.Item($Row, 1)
    $dateCell.Value = $dateCell.Value2() + 60

    $idCell = $ws.Range($IdCol + $Row)
    $parts = $idCell.Value().Split("|")
    $oldDate = [DateTime]::ParseExact($parts[0], "yyyy-MM-dd", $null)
    $parts[0] = $oldDate.AddDays(60).ToString("yyyy-MM-dd")
    $idCell.Value = [string]::Join("|", $parts)
}

Shift-DailyDetailsRow 2 "V"
Shift-DailyDetailsRow 3 "V"
Shift-DailyDetailsRow 4 "V"
Shift-DailyDetailsRow 5 "V"
Shift-DailyDetailsRow 6 "V"

Shift-DailyDetailsRow 7 "W"
Shift-DailyDetailsRow 8 "W"
Shift-DailyDetailsRow 9 "W"

Shift-DailyDetailsRow 10 "W"
Shift-DailyDetailsRow 11 "W"
Shift-DailyDetailsRow 12 "W"
Shift-DailyDetailsRow 13 "W"
Shift-DailyDetailsRow 14 "W"
Shift-DailyDetailsRow 15 "W"
Shift-DailyDetailsRow 16 "W"
